$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-6: 45207 -> 45208
foreach ($r in 2..6) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
